$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new data row (row 16), following the same layout as the
# preceding row 15: col A = running index, col B = category label,
# cols C:M = intensity values of 1.

# Copy formatting (incl. style) from A15 -> A16, then set the new value.
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null
$ws.Range("A16").Value2 = 14

# B16 reuses the same text label as B15 ("HexGrid-60degTilt5degRes").
$ws.Range("B16").Value2 = $ws.Range("B15").Value2

# C16:M16 are all 1.
$ws.Range("C16:M16").Value2 = 1

$excel.CutCopyMode = 0
